# Applies the diff: updates Brasil/Nordeste/Sergipe 'Roubo seguido de morte' data
# by inserting 2025 rows for each region and correcting recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.Style = 'Normal'
}

function Set-NewDataRow($ws, $rowNum, $region, $date, $variable, $valor, $posicao, $faltam) {
    Set-TextCell $ws "A$rowNum" $region
    Set-TextCell $ws "B$rowNum" $date
    Set-TextCell $ws "C$rowNum" $variable
    $ws.Range("D$rowNum").Value = $valor
    if ($null -ne $posicao) {
        $ws.Range("E$rowNum").Value = $posicao
    }
    $ws.Range("F$rowNum").Value = $faltam
}

# 1) Insert new row for Brasil 01/01/2025 at the end of the Brasil block (row 12)
$ws.Rows.Item(12).Insert()

# 2) Insert new row for Nordeste 01/01/2025 at the end of the (now shifted) Nordeste block (row 23)
$ws.Rows.Item(23).Insert()

Set-NewDataRow $ws 12 'Brasil' '01/01/2025' 'Roubo seguido de morte (latrocínio)' 0.2265864023064676 $null $false
Set-NewDataRow $ws 23 'Nordeste' '01/01/2025' 'Roubo seguido de morte (latrocínio)' 0.243395964067396 $null $false

# 3) Append new row for Sergipe 01/01/2025 directly after the existing last Sergipe row (row 34)
Set-NewDataRow $ws 34 'Sergipe' '01/01/2025' 'Roubo seguido de morte (latrocínio)' 0.1243821317604798 24 $false

# 4) Apply the recalculated values for rows whose underlying statistics changed
$ws.Range('D3').Value = 1.484195478947663
$ws.Range('D4').Value = 1.491599847558214
$ws.Range('D5').Value = 1.292147025239075
$ws.Range('D7').Value = 1.014607599515074
$ws.Range('D9').Value = 0.7708129718437964
$ws.Range('D10').Value = 0.5835148157186814
$ws.Range('D11').Value = 0.4747420516039645
$ws.Range('D20').Value = 0.7944406142306738
$ws.Range('D21').Value = 0.579551901079609
$ws.Range('D22').Value = 0.5070466835479671
$ws.Range('E29').Value = 9
